$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.82'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.43%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.70'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.45%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.049'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.16%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07993'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.09%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.865'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-4.34%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.773'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.06%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9214'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.79%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1290'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-5.49%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1882'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.29%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09096'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.39%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03422'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.77%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09879'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.01%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001415'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.41%'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '6.59%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.855'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '8.13%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.115'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.99%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.394'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '14.04%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3417'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.84%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1339'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.42%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.01%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2497'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.70%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04413'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.16%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001233'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.03%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004882'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.32%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001300'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-21.22%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '42.18%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01941'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-1.46%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05172'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.56%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007564'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.87%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01012'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-7.74%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1353'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.02%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002131'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.48%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009907'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-8.29%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006187'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-3.46%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.05%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.01'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-0.33%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001250'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '5.03%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.05%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.05%'
